$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.276739466669731
$ws.Range("C2").Value = 5.512395282490526
$ws.Range("D2").Value = 4.922777254636135
$ws.Range("F2").Value = 22.84158384033515
$ws.Range("G2").Value = 3.626838973828984
$ws.Range("K2").Value = 8.50592536256166
$ws.Range("N2").Value = 17.80671404251862
$ws.Range("O2").Value = 20.5876497060342
$ws.Range("B3").Value = 8.976158137798656
$ws.Range("C3").Value = 5.385868133065467
$ws.Range("D3").Value = 4.838897443180038
$ws.Range("F3").Value = 22.9019181309469
$ws.Range("G3").Value = 3.628319191498725
$ws.Range("K3").Value = 8.296656281057613
$ws.Range("N3").Value = 17.85861684992159
$ws.Range("O3").Value = 20.67349291331158
$ws.Range("B4").Value = 8.78771250488265
$ws.Range("C4").Value = 5.306062920337981
$ws.Range("D4").Value = 4.785988330907564
$ws.Range("F4").Value = 22.94512450720505
$ws.Range("G4").Value = 3.629276656475332
$ws.Range("K4").Value = 8.166496745428491
$ws.Range("N4").Value = 17.8921649449647
$ws.Range("O4").Value = 20.73057661869371
$ws.Range("B5").Value = 8.71006266017603
$ws.Range("C5").Value = 5.273039829396891
$ws.Range("D5").Value = 4.76408841237421
$ws.Range("F5").Value = 22.96427654796425
$ws.Range("G5").Value = 3.62967908976583
$ws.Range("K5").Value = 8.113112369215932
$ws.Range("N5").Value = 17.90625941430286
$ws.Range("O5").Value = 20.75493746427076
$ws.Range("B6").Value = 8.697120939835996
$ws.Range("C6").Value = 5.267526988988356
$ws.Range("D6").Value = 4.76043191159664
$ws.Range("F6").Value = 22.96754992138374
$ws.Range("G6").Value = 3.62974665496266
$ws.Range("K6").Value = 8.104229490659726
$ws.Range("N6").Value = 17.90862538920466
$ws.Range("O6").Value = 20.75904888292994
$ws.Range("B7").Value = 8.786668593964984
$ws.Range("C7").Value = 5.305619549577028
$ws.Range("D7").Value = 4.785694333840927
$ws.Range("F7").Value = 22.94537654803796
$ws.Range("G7").Value = 3.629282034148324
$ws.Range("K7").Value = 8.165778072038512
$ws.Range("N7").Value = 17.89235331244148
$ws.Range("O7").Value = 20.73090071069399
$ws.Range("B8").Value = 9.173975338241547
$ws.Range("C8").Value = 5.469225860716903
$ws.Range("D8").Value = 4.894156036768734
$ws.Range("F8").Value = 22.86110606162387
$ws.Range("G8").Value = 3.627339286746626
$ws.Range("K8").Value = 8.434159906572702
$ws.Range("N8").Value = 17.82426214534002
$ws.Range("O8").Value = 20.61633937532591
$ws.Range("B9").Value = 9.897826890339617
$ws.Range("C9").Value = 5.772012871345652
$ws.Range("D9").Value = 5.095077272935579
$ws.Range("F9").Value = 22.74490257261192
$ws.Range("G9").Value = 3.623913566350925
$ws.Range("K9").Value = 8.944082530455237
$ws.Range("N9").Value = 17.70401586947312
$ws.Range("O9").Value = 20.42647087089439
$ws.Range("B10").Value = 10.40215545890946
$ws.Range("C10").Value = 5.981934960164851
$ws.Range("D10").Value = 5.234747377261789
$ws.Range("F10").Value = 22.6896328749126
$ws.Range("G10").Value = 3.621628459990907
$ws.Range("K10").Value = 9.304830990261683
$ws.Range("N10").Value = 17.62370186199631
$ws.Range("O10").Value = 20.30827221892187
$ws.Range("B11").Value = 10.62459616917833
$ws.Range("C11").Value = 6.074425196451371
$ws.Range("D11").Value = 5.296412021854102
$ws.Range("F11").Value = 22.67105986141523
$ws.Range("G11").Value = 3.620638738034809
$ws.Range("K11").Value = 9.465191181101089
$ws.Range("N11").Value = 17.58889506459125
$ws.Range("O11").Value = 20.25914472662654
$ws.Range("B12").Value = 10.70775170678241
$ws.Range("C12").Value = 6.108995338098506
$ws.Range("D12").Value = 5.319482029691986
$ws.Range("F12").Value = 22.66497331486777
$ws.Range("G12").Value = 3.62027107772369
$ws.Range("K12").Value = 9.525321537435286
$ws.Range("N12").Value = 17.57596218517839
$ws.Range("O12").Value = 20.24121055875612
$ws.Range("B13").Value = 10.68989170600194
$ws.Range("C13").Value = 6.101570555999145
$ws.Range("D13").Value = 5.314526171958689
$ws.Range("F13").Value = 22.66624203345842
$ws.Range("G13").Value = 3.62034994354544
$ws.Range("K13").Value = 9.512398673579812
$ws.Range("N13").Value = 17.57873650902146
$ws.Range("O13").Value = 20.24504320351015
$ws.Range("B14").Value = 10.63145942324439
$ws.Range("C14").Value = 6.077278518539626
$ws.Range("D14").Value = 5.298315703435778
$ws.Range("F14").Value = 22.67054013592787
$ws.Range("G14").Value = 3.620608347773889
$ws.Range("K14").Value = 9.470150361393612
$ws.Range("N14").Value = 17.58782610936329
$ws.Range("O14").Value = 20.25765584640749
$ws.Range("B15").Value = 10.59552552890998
$ws.Range("C15").Value = 6.062339214428598
$ws.Range("D15").Value = 5.28834938195084
$ws.Range("F15").Value = 22.67329617884922
$ws.Range("G15").Value = 3.620767554853765
$ws.Range("K15").Value = 9.444193040423499
$ws.Range("N15").Value = 17.59342598485856
$ws.Range("O15").Value = 20.26546867869657
$ws.Range("B16").Value = 10.38747143525045
$ws.Range("C16").Value = 5.975828137905811
$ws.Range("D16").Value = 5.230678747785558
$ws.Range("F16").Value = 22.6909790401051
$ws.Range("G16").Value = 3.621694139608154
$ws.Range("K16").Value = 9.29427081072388
$ws.Range("N16").Value = 17.62601127589091
$ws.Range("O16").Value = 20.3115763794501
$ws.Range("B17").Value = 10.25799359119706
$ws.Range("C17").Value = 5.921970751400464
$ws.Range("D17").Value = 5.194811549600097
$ws.Range("F17").Value = 22.70351107118528
$ws.Range("G17").Value = 3.62227529673211
$ws.Range("K17").Value = 9.201297678992439
$ws.Range("N17").Value = 17.64644338412863
$ws.Range("O17").Value = 20.34105214599798
$ws.Range("B18").Value = 10.18286815609409
$ws.Range("C18").Value = 5.890711810696005
$ws.Range("D18").Value = 5.174006149269069
$ws.Range("F18").Value = 22.71133742768065
$ws.Range("G18").Value = 3.622614250836883
$ws.Range("K18").Value = 9.14747243714527
$ws.Range("N18").Value = 17.65835813180871
$ws.Range("O18").Value = 20.35844261138046
$ws.Range("B19").Value = 10.15732218309982
$ws.Range("C19").Value = 5.880080402749843
$ws.Range("D19").Value = 5.166932022998429
$ws.Range("F19").Value = 22.71409341593807
$ws.Range("G19").Value = 3.622729821034682
$ws.Range("K19").Value = 9.129189877418629
$ws.Range("N19").Value = 17.66242023744176
$ws.Range("O19").Value = 20.36440568467864
$ws.Range("B20").Value = 10.2718448963473
$ws.Range("C20").Value = 5.927733270601532
$ws.Range("D20").Value = 5.198647929687934
$ws.Range("F20").Value = 22.70211301110319
$ws.Range("G20").Value = 3.622212946624385
$ws.Range("K20").Value = 9.211231424769396
$ws.Range("N20").Value = 17.6442515146046
$ws.Range("O20").Value = 20.3378691811394
$ws.Range("B21").Value = 10.64865220476612
$ws.Range("C21").Value = 6.084426162905912
$ws.Range("D21").Value = 5.303084829713344
$ws.Range("F21").Value = 22.66925197421645
$ws.Range("G21").Value = 3.620532255054286
$ws.Range("K21").Value = 9.482576273641019
$ws.Range("N21").Value = 17.58514955760938
$ws.Range("O21").Value = 20.25393302664897
$ws.Range("B22").Value = 10.88860597723619
$ws.Range("C22").Value = 6.184179455355665
$ws.Range("D22").Value = 5.369697183776168
$ws.Range("F22").Value = 22.65329371657098
$ws.Range("G22").Value = 3.619475346121021
$ws.Range("K22").Value = 9.656427540186286
$ws.Range("N22").Value = 17.54796641053065
$ws.Range("O22").Value = 20.20297822647332
$ws.Range("B23").Value = 10.76113798113192
$ws.Range("C23").Value = 6.131188898236926
$ws.Range("D23").Value = 5.334299002923815
$ws.Range("F23").Value = 22.66130550201726
$ws.Range("G23").Value = 3.620035650046147
$ws.Range("K23").Value = 9.563976297849857
$ws.Range("N23").Value = 17.56767996142618
$ws.Range("O23").Value = 20.22981608168759
$ws.Range("B24").Value = 10.26558485510212
$ws.Range("C24").Value = 5.925128954772214
$ws.Range("D24").Value = 5.196914077175235
$ws.Range("F24").Value = 22.70274313823254
$ws.Range("G24").Value = 3.622241120032189
$ws.Range("K24").Value = 9.206741538150821
$ws.Range("N24").Value = 17.64524193526449
$ws.Range("O24").Value = 20.33930681525108
$ws.Range("B25").Value = 9.706454375895301
$ws.Range("C25").Value = 5.692198645300256
$ws.Range("D25").Value = 5.0420651878804
$ws.Range("F25").Value = 22.77106471649386
$ws.Range("G25").Value = 3.624799446122124
$ws.Range("K25").Value = 8.80830885621827
$ws.Range("N25").Value = 17.73513071227949
$ws.Range("O25").Value = 20.47410180718193
